$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.862.16'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.629.99'
$ws.Range("D3").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5067'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2576'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06327'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.64%  '
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07754'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.246'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.634.21'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.853.78'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5502'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.66'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₅7651'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.877.97'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.33'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.409'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.871'
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.024'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.912'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.99'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1245'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.34%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.60'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.765'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.239'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04887'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.244'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.190'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.544'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.368'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8951'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5526'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.537'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.119.25'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01553'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.584'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7960'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.33'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₈119'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.765.67'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4443'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.73'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05135'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.545'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.92%  '
